# liensAllUserFR.xlsx – "Add files via upload" edit
#
# Content change: every section that previously showed the placeholder
# "Rien pour le moment…" (with an empty link cell) is rewritten to the new
# 4-column placeholder layout used elsewhere in the sheet:
#   B = "Liens : "   (trailing space, style copied from A1 -> no vertical-center)
#   C = "Rien pour le moment"   (no ellipsis)
#   D = "rien"
#   E = "ignore"     (marker column, default/no style)
# This touches rows 23, 47, 48, 49, 50, 51, 52, 53 and 59.
# The selection is also left on B23:E23, matching the author's last edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 already carries the target "no special alignment" style (cellXfs #1:
# Arial 10, no vertical-center) that B/C/D need to end up with, so grab the
# format from there instead of re-building fonts/alignment by hand (which
# would otherwise mint new, slightly different style/font entries).
$ws.Range("A1").Copy()

$rows = 23, 47, 48, 49, 50, 51, 52, 53, 59

foreach ($r in $rows) {
    $target = "B" + $r + ":D" + $r
    [void]$ws.Range($target).PasteSpecial(-4122)

    $ws.Range("B" + $r).Value = "Liens : "
    $ws.Range("C" + $r).Value = "Rien pour le moment"
    $ws.Range("D" + $r).Value = "rien"
    $ws.Range("E" + $r).Value = "ignore"
}

$excel.CutCopyMode = 0

[void]$ws.Range("B23:E23").Select()
